# pendulum experiment new added
# Update the data values in rows 2-4 (columns B:I) with new experiment data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 4.4
$ws.Range("C2").Value = 5.5
$ws.Range("D2").Value = 1.3
$ws.Range("E2").Value = 2.41
$ws.Range("F2").Value = 1.1
$ws.Range("G2").Value = -0.6172624830350726
$ws.Range("H2").Value = -0.09824037536020899
$ws.Range("I2").Value = 5.711986642890535

# Row 3
$ws.Range("B3").Value = 2.180800269
$ws.Range("C3").Value = 3.881050218
$ws.Range("D3").Value = 2.412134171
$ws.Range("E3").Value = 0.2406511307
$ws.Range("F3").Value = 1.700249949
$ws.Range("G3").Value = 2.304918888488151
$ws.Range("H3").Value = 0.3668392345287664
$ws.Range("I3").Value = 3.6954480197897

# Row 4
$ws.Range("B4").Value = 5.081199216
$ws.Range("C4").Value = 6.880759618
$ws.Range("D4").Value = 0.2406511307
$ws.Range("E4").Value = 0.05192375183
$ws.Range("F4").Value = 1.799560402
$ws.Range("G4").Value = 1.533571962728984
$ws.Range("H4").Value = 0.2440755584554578
$ws.Range("I4").Value = 3.491511204734536
